$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 382, shifting existing rows 382:485 down to 383:486.
$ws.Rows.Item(382).Insert()

# Populate the newly inserted row 382 with the new record's data.
$ws.Cells.Item(382, 1).Value = 3
$ws.Cells.Item(382, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(382, 3).Value = "Coquimbo"
$ws.Cells.Item(382, 4).Value = 44736
$ws.Cells.Item(382, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(382, 5).Value = 5
$ws.Cells.Item(382, 6).Value = 100112032
$ws.Cells.Item(382, 7).Value = "Zapallo italiano"
$ws.Cells.Item(382, 8).Value = "Sin especificar"
$ws.Cells.Item(382, 9).Value = "Primera"
$ws.Cells.Item(382, 10).Value = 165
$ws.Cells.Item(382, 11).Value = 9000
$ws.Cells.Item(382, 12).Value = 10000
$ws.Cells.Item(382, 13).Value = 9515
$ws.Cells.Item(382, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(382, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(382, 16).Value = 136
$ws.Cells.Item(382, 17).Value = 70
$ws.Cells.Item(382, 18).Value = "Hortaliza"
